$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: insert "Yes"/"No" rows before "Skip", and turn old "ComfirmBox" row into "Exit" ---
$tmp = $ws2.Rows.Item(5).Insert()
$tmp = $ws2.Rows.Item(5).Insert()
$ws2.Range("A5").Value = "Yes"
$ws2.Range("A6").Value = "No"
$ws2.Range("A9").Value = "Exit"

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- Sheet1: drop the old "ShouldStopRead" column so column headers become Index/Type/GoTo/Script ---
$ws1.Columns.Item(2).Delete()

# add a sample data row
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "Exit"

# --- data validation lists driven from Sheet2 ---
$ws1.Range("B34").Validation.Add(3, 1, 1, 'Sheet2!$A$2:$A$8')
$ws1.Range("B2:B33").Validation.Add(3, 1, 1, 'Sheet2!$A$2:$A$10')

# --- view/selection state ---
$tmp = $ws2.Range("A12").Select()
$tmp = $ws1.Activate()
$tmp = $ws1.Range("C2").Select()
